# legado: PMQA - PONTOS
# Rename header labels in row 1 to the new snake_case naming convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "nome_ponto_coleta"
$ws.Range("I1").Value = "bacia_hidrografica"
$ws.Range("F1").Value = "Tipo_ambiente"

# Update the active selection to match the saved view state.
$ws.Range("F5").Select()
